$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Tests for shorts: rework observed readings for two connections ---
$ws.Range("E14").Value = 20
$ws.Range("E15").Value = 10

# --- Flex cable continuity checks reworked to "infinite" (open) reading ---
$ws.Range("E27").Value = "∞"

# --- HV/leakage style readings tightened slightly ---
$ws.Range("E30").Value = 90.5
$ws.Range("E32").Value = 91

# --- More re-tested adapters: shorts rework + infinite reading retest ---
$ws.Range("E40").Value = 30
$ws.Range("E41").Value = "∞"
$ws.Range("E42").Value = "∞"
$ws.Range("E43").Value = "∞"

# --- Restore the on-screen view/selection state recorded at save time ---
$ws.Range("E44").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1

$excel.Calculate()
